# Updates the cryptos worksheet: Price (D) and Volume(1h) (E) columns for rows 2-51.
# Mirrors a refreshed crypto-price snapshot (GitHub Actions scheduled run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.506.05"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.570.98"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -1.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.14%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.795.37"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "1.578.59"
$ws.Range("E13").Value = "  +0.65%  "
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "27.468.75"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "1.455.02"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0168"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  +6.57%  "
$ws.Range("E45").Value = "  -3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.28%  "
$ws.Range("D47").Value = "1.707.28"
$ws.Range("E47").Value = "  +0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0525"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.45%  "
$ws.Range("E51").Value = "  -1.57%  "
